$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Ithaca Soda - Ginger Beer - quantity 1 -> 2, total cost 28.75 -> 57.50
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "57.50"

# Row 5: Ithaca Soda - Root Beer - quantity 1 -> 2, total cost 28.75 -> 57.50
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "2"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "57.50"
